# EEE2000-metadata_mapping.xlsx — data fixes
# 1) "hasPart" -> "HasPart" (casing fix) in the CCI/Aerosol relationship cells
# 2) Malformed date-like text values (day/month transposed, e.g. "2022-31-12")
#    replaced with real Excel dates

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

# --- Fix casing: hasPart, -> HasPart, (rows 101-109, column G) ---
$ws.Range("G101:G109").Value = "HasPart,`nIsNewVersionOf"

# --- Fix bogus date text "2022-31-12" -> real date 2022-12-31 (rows 128-141, column C) ---
$ws.Range("C128:C141").Value = "12/31/2022"

# --- Fix bogus date text "2021-09-31" -> real date 2021-09-30 (row 147, column C) ---
$ws.Range("C147").Value = "9/30/2021"

# --- Fix bogus date text "2018-09-31" -> real date 2018-09-30 (row 148, column C) ---
$ws.Range("C148").Value = "9/30/2018"
